$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 195, shifting existing rows 195:244 down to 196:245.
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new record's values.
$ws.Range("A195").Value = 4
$ws.Range("B195").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C195").Value = "Los Lagos"
$ws.Range("D195").Value = 44642
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E195").Value = 10
$ws.Range("F195").Value = 100112043
$ws.Range("G195").Value = "Pepino ensalada"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 400
$ws.Range("K195").Value = 21000
$ws.Range("L195").Value = 21000
$ws.Range("M195").Value = 21000
$ws.Range("N195").Value = "$/caja 60 unidades"
$ws.Range("O195").Value = "Región de Arica y Parinacota"
$ws.Range("P195").Value = 350
$ws.Range("Q195").Value = 60
$ws.Range("R195").Value = "Hortaliza"
